$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the home/away team names in row 2:
#  - A2 ("Varberg") becomes "Kungsbacka"
#  - B2 ("Kungsbacka") becomes "Trollhättan"
$ws.Range("A2").Value = "Kungsbacka"
$ws.Range("B2").Value = "Trollhättan"

# B2 should pick up the same direct formatting as B1 (style index 2:
# Arial/black font, solid white fill, left-aligned) instead of the
# generic style used previously.
$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# Remove the "Marten Gullberg" player row entirely (row 6).
$ws.Rows(6).Delete()

Write-Host "done"
